# Update column G ("K") values on Sheet1 for rows 2-28.
# This mirrors the upstream regeneration of save_data that now writes the
# strike-count column ("K") using the newly calculated s_vals instead of
# the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 1
    6  = 2
    7  = 0
    8  = 2
    9  = 1
    10 = 2
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 2
    17 = 0
    18 = 2
    19 = 1
    20 = 2
    21 = 1
    22 = 2
    23 = 0
    24 = 2
    25 = 1
    26 = 1
    27 = 1
    28 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
